$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.624526
$ws.Range("H2").Value = 7.873578
$ws.Range("I2").Value = 0.06442870872595916
$ws.Range("J2").Value = 0.06442870872595916
$ws.Range("M2").Value = 47.32925566666668
$ws.Range("N2").Value = 141.987767
$ws.Range("O2").Value = 0.3408416299313156
$ws.Range("P2").Value = 0.3408416299313156
$ws.Range("Q2").Value = 124.216862057814
$ws.Range("R2").Value = 1117.951758520326
$ws.Range("S2").Value = 0.02195998609652589
$ws.Range("T2").Value = 0.0219599860965259
$ws.Range("G3").Value = 2.624526
$ws.Range("H3").Value = 7.873578
$ws.Range("I3").Value = 0.06442870872595916
$ws.Range("J3").Value = 0.06442870872595916
$ws.Range("M3").Value = 43.717953
$ws.Range("O3").Value = 0.3148348341399153
$ws.Range("P3").Value = 0.3148348341399154
$ws.Range("Q3").Value = 114.738904315278
$ws.Range("R3").Value = 1032.650138837502
$ws.Range("S3").Value = 0.02028440182558626
$ws.Range("T3").Value = 0.02028440182558627
$ws.Range("G4").Value = 2.624526
$ws.Range("H4").Value = 7.873578
$ws.Range("I4").Value = 0.06442870872595916
$ws.Range("J4").Value = 0.06442870872595916
$ws.Range("M4").Value = 21.09134933333333
$ws.Range("N4").Value = 63.274048
$ws.Range("O4").Value = 0.1518893501062827
$ws.Range("P4").Value = 0.1518893501062827
$ws.Range("Q4").Value = 55.354794700416
$ws.Range("R4").Value = 498.193152303744
$ws.Range("S4").Value = 0.009786034696572922
$ws.Range("T4").Value = 0.009786034696572923
$ws.Range("G5").Value = 2.624526
$ws.Range("H5").Value = 7.873578
$ws.Range("I5").Value = 0.06442870872595916
$ws.Range("J5").Value = 0.06442870872595916
$ws.Range("M5").Value = 26.72140366666666
$ws.Range("N5").Value = 80.16421099999999
$ws.Range("O5").Value = 0.1924341858224863
$ws.Range("P5").Value = 0.1924341858224864
$ws.Range("Q5").Value = 70.13101867966199
$ws.Range("R5").Value = 631.179168116958
$ws.Range("S5").Value = 0.01239828610727407
$ws.Range("T5").Value = 0.01239828610727407
$ws.Range("I6").Value = 0.01210207677934112
$ws.Range("J6").Value = 0.01210207677934112
$ws.Range("M6").Value = 47.32925566666668
$ws.Range("N6").Value = 141.987767
$ws.Range("O6").Value = 0.3408416299313156
$ws.Range("P6").Value = 0.3408416299313156
$ws.Range("Q6").Value = 23.33248689348323
$ws.Range("R6").Value = 209.992382041349
$ws.Range("S6").Value = 0.004124891575024554
$ws.Range("T6").Value = 0.004124891575024555
$ws.Range("I7").Value = 0.01210207677934112
$ws.Range("J7").Value = 0.01210207677934112
$ws.Range("M7").Value = 43.717953
$ws.Range("O7").Value = 0.3148348341399153
$ws.Range("P7").Value = 0.3148348341399154
$ws.Range("S7").Value = 0.003810155335572382
$ws.Range("T7").Value = 0.003810155335572384
$ws.Range("I8").Value = 0.01210207677934112
$ws.Range("J8").Value = 0.01210207677934112
$ws.Range("M8").Value = 21.09134933333333
$ws.Range("N8").Value = 63.274048
$ws.Range("O8").Value = 0.1518893501062827
$ws.Range("P8").Value = 0.1518893501062827
$ws.Range("Q8").Value = 10.39766260749511
$ws.Range("R8").Value = 93.57896346745601
$ws.Range("S8").Value = 0.001838176576950458
$ws.Range("T8").Value = 0.001838176576950459
$ws.Range("I9").Value = 0.01210207677934112
$ws.Range("J9").Value = 0.01210207677934112
$ws.Range("M9").Value = 26.72140366666666
$ws.Range("N9").Value = 80.16421099999999
$ws.Range("O9").Value = 0.1924341858224863
$ws.Range("P9").Value = 0.1924341858224864
$ws.Range("Q9").Value = 13.17317992953522
$ws.Range("R9").Value = 118.558619365817
$ws.Range("S9").Value = 0.002328853291793726
$ws.Range("T9").Value = 0.002328853291793727
$ws.Range("I10").Value = 0.003429134645952472
$ws.Range("J10").Value = 0.003429134645952472
$ws.Range("M10").Value = 47.32925566666668
$ws.Range("N10").Value = 141.987767
$ws.Range("O10").Value = 0.3408416299313156
$ws.Range("P10").Value = 0.3408416299313156
$ws.Range("Q10").Value = 6.611281736309668
$ws.Range("R10").Value = 59.50153562678701
$ws.Range("S10").Value = 0.001168791841980385
$ws.Range("T10").Value = 0.001168791841980385
$ws.Range("I11").Value = 0.003429134645952472
$ws.Range("J11").Value = 0.003429134645952472
$ws.Range("M11").Value = 43.717953
$ws.Range("O11").Value = 0.3148348341399153
$ws.Range("P11").Value = 0.3148348341399154
$ws.Range("Q11").Value = 6.106829700711001
$ws.Range("R11").Value = 54.96146730639901
$ws.Range("S11").Value = 0.001079611037501883
$ws.Range("T11").Value = 0.001079611037501884
$ws.Range("I12").Value = 0.003429134645952472
$ws.Range("J12").Value = 0.003429134645952472
$ws.Range("M12").Value = 21.09134933333333
$ws.Range("N12").Value = 63.274048
$ws.Range("O12").Value = 0.1518893501062827
$ws.Range("P12").Value = 0.1518893501062827
$ws.Range("Q12").Value = 2.946187314325333
$ws.Range("R12").Value = 26.515685828928
$ws.Range("S12").Value = 0.0005208490328006588
$ws.Range("T12").Value = 0.0005208490328006589
$ws.Range("I13").Value = 0.003429134645952472
$ws.Range("J13").Value = 0.003429134645952472
$ws.Range("M13").Value = 26.72140366666666
$ws.Range("N13").Value = 80.16421099999999
$ws.Range("O13").Value = 0.1924341858224863
$ws.Range("P13").Value = 0.1924341858224864
$ws.Range("Q13").Value = 3.732632713985666
$ws.Range("R13").Value = 33.593694425871
$ws.Range("S13").Value = 0.0006598827336695437
$ws.Range("T13").Value = 0.0006598827336695439
$ws.Range("G14").Value = 37.47815466666666
$ws.Range("H14").Value = 112.434464
$ws.Range("I14").Value = 0.9200400798487472
$ws.Range("J14").Value = 0.9200400798487472
$ws.Range("M14").Value = 47.32925566666668
$ws.Range("N14").Value = 141.987767
$ws.Range("O14").Value = 0.3408416299313156
$ws.Range("P14").Value = 0.3408416299313156
$ws.Range("Q14").Value = 1773.813164133543
$ws.Range("R14").Value = 15964.31847720189
$ws.Range("S14").Value = 0.3135879604177848
$ws.Range("T14").Value = 0.3135879604177848
$ws.Range("G15").Value = 37.47815466666666
$ws.Range("H15").Value = 112.434464
$ws.Range("I15").Value = 0.9200400798487472
$ws.Range("J15").Value = 0.9200400798487472
$ws.Range("M15").Value = 43.717953
$ws.Range("O15").Value = 0.3148348341399153
$ws.Range("P15").Value = 0.3148348341399154
$ws.Range("Q15").Value = 1638.468204244064
$ws.Range("R15").Value = 14746.21383819658
$ws.Range("S15").Value = 0.2896606659412547
$ws.Range("T15").Value = 0.2896606659412548
$ws.Range("G16").Value = 37.47815466666666
$ws.Range("H16").Value = 112.434464
$ws.Range("I16").Value = 0.9200400798487472
$ws.Range("J16").Value = 0.9200400798487472
$ws.Range("M16").Value = 21.09134933333333
$ws.Range("N16").Value = 63.274048
$ws.Range("O16").Value = 0.1518893501062827
$ws.Range("P16").Value = 0.1518893501062827
$ws.Range("Q16").Value = 790.4648524433635
$ws.Range("R16").Value = 7114.183671990271
$ws.Range("S16").Value = 0.1397442897999587
$ws.Range("T16").Value = 0.1397442897999587
$ws.Range("G17").Value = 37.47815466666666
$ws.Range("H17").Value = 112.434464
$ws.Range("I17").Value = 0.9200400798487472
$ws.Range("J17").Value = 0.9200400798487472
$ws.Range("M17").Value = 26.72140366666666
$ws.Range("N17").Value = 80.16421099999999
$ws.Range("O17").Value = 0.1924341858224863
$ws.Range("P17").Value = 0.1924341858224864
$ws.Range("Q17").Value = 1001.468899529767
$ws.Range("R17").Value = 9013.220095767903
$ws.Range("S17").Value = 0.177047163689749
$ws.Range("T17").Value = 0.177047163689749
